# Update Result/Date columns (A/B) for each data sheet, marking rows as
# Fail and stamping the new run timestamp, as produced by the Katalon
# "Demo Verification Script" run after the Windows 11 update.

$wb = $excel.ActiveWorkbook

$wsSale = $wb.Worksheets.Item("CC-Payments-Sale")
$wsSale.Range("A2").Value = "Fail"
$wsSale.Range("B2").Value = "Mon Aug 04 20:40:37 IST 2025"

$wsAuth = $wb.Worksheets.Item("CC-Payments-Auth")
$wsAuth.Range("A2").Value = "Fail"
$wsAuth.Range("B2").Value = "Mon Aug 04 20:27:19 IST 2025"
$wsAuth.Range("A3").Value = "Fail"
$wsAuth.Range("B3").Value = "Mon Aug 04 20:28:16 IST 2025"
$wsAuth.Range("A4").Value = "Fail"
$wsAuth.Range("B4").Value = "Mon Aug 04 20:28:55 IST 2025"
$wsAuth.Range("A5").Value = "Fail"
$wsAuth.Range("B5").Value = "Mon Aug 04 20:29:41 IST 2025"
$wsAuth.Range("A6").Value = "Fail"
$wsAuth.Range("B6").Value = "Mon Aug 04 20:30:20 IST 2025"
$wsAuth.Range("A7").Value = "Fail"
$wsAuth.Range("B7").Value = "Mon Aug 04 20:31:07 IST 2025"

$wsDebit = $wb.Worksheets.Item("ACH-Payments-Debit")
$wsDebit.Range("A2").Value = "Fail"
$wsDebit.Range("B2").Value = "Mon Aug 04 20:32:20 IST 2025"
$wsDebit.Range("A3").Value = "Fail"
$wsDebit.Range("B3").Value = "Mon Aug 04 20:33:08 IST 2025"
$wsDebit.Range("A4").Value = "Fail"
$wsDebit.Range("B4").Value = "Mon Aug 04 20:33:49 IST 2025"
$wsDebit.Range("A5").Value = "Fail"
$wsDebit.Range("B5").Value = "Mon Aug 04 20:34:38 IST 2025"
$wsDebit.Range("A6").Value = "Fail"
$wsDebit.Range("B6").Value = "Mon Aug 04 20:35:25 IST 2025"
$wsDebit.Range("A7").Value = "Fail"
$wsDebit.Range("B7").Value = "Mon Aug 04 20:36:13 IST 2025"
$wsDebit.Range("A8").Value = "Fail"
$wsDebit.Range("B8").Value = "Mon Aug 04 20:37:01 IST 2025"
$wsDebit.Range("A9").Value = "Fail"
$wsDebit.Range("B9").Value = "Mon Aug 04 20:37:47 IST 2025"
$wsDebit.Range("A10").Value = "Fail"
$wsDebit.Range("B10").Value = "Mon Aug 04 20:38:41 IST 2025"
